$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1870
$ws.Range("J43").Value = 1740
$ws.Range("L43").Value = 1740
$ws.Range("N43").Value = -1878
$ws.Range("H70").Value = 1149.7858
$ws.Range("I70").Value = 1122.8462
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 3368.5386
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -3098.5386
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 1149.7858
$ws.Range("I73").Value = 1122.8462
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 3368.5386
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -2432.5386
$ws.Range("N73").Value = -6372
$ws.Range("H86").Value = 8456.357
$ws.Range("I86").Value = 1467.5714
$ws.Range("J86").Value = 15445.143
$ws.Range("K86").Value = 1467.5714
$ws.Range("L86").Value = 15445.143
$ws.Range("M86").Value = -344.5714
$ws.Range("N86").Value = -17691.143
$ws.Range("H89").Value = 8456.357
$ws.Range("I89").Value = 1467.5714
$ws.Range("J89").Value = 15445.143
$ws.Range("K89").Value = 7337.857
$ws.Range("L89").Value = 77225.715
$ws.Range("M89").Value = -1721.857
$ws.Range("N89").Value = -88457.715
$ws.Range("H129").Value = 836.4394
$ws.Range("J129").Value = 841.3492
$ws.Range("L129").Value = 2524.0476
$ws.Range("N129").Value = -12524.0476
$ws.Range("H132").Value = 3098.28
$ws.Range("I132").Value = 3343.2273
$ws.Range("K132").Value = 10029.6819
$ws.Range("M132").Value = -7499.6819
$ws.Range("H137").Value = 2070.4348
$ws.Range("I137").Value = 2051.875
$ws.Range("J137").Value = 2112.8572
$ws.Range("K137").Value = 6155.625
$ws.Range("L137").Value = 6338.571599999999
$ws.Range("M137").Value = -3605.625
$ws.Range("N137").Value = -11438.5716
$ws.Range("H138").Value = 1927.825
$ws.Range("I138").Value = 531.7692
$ws.Range("J138").Value = 4520.5
$ws.Range("K138").Value = 1595.3076
$ws.Range("L138").Value = 13561.5
$ws.Range("M138").Value = 3544.6924
$ws.Range("N138").Value = -23841.5
$ws.Range("H141").Value = 4084.2856
$ws.Range("I141").Value = 3672.5
$ws.Range("K141").Value = 11017.5
$ws.Range("M141").Value = -5837.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2120.578
$ws.Range("I32").Value = 1836.7632
$ws.Range("J32").Value = 3661.2856
$ws.Range("K32").Value = 1836.7632
$ws.Range("L32").Value = 3661.2856
$ws.Range("M32").Value = -1549.7632
$ws.Range("N32").Value = -4235.2856
$ws.Range("H74").Value = 2390.0557
$ws.Range("I74").Value = 2089.4707
$ws.Range("K74").Value = 2089.4707
$ws.Range("M74").Value = -1215.4707
$ws.Range("H77").Value = 2390.0557
$ws.Range("I77").Value = 2089.4707
$ws.Range("K77").Value = 10447.3535
$ws.Range("M77").Value = -6079.353499999999
$ws.Range("H102").Value = 3229.25
$ws.Range("I102").Value = 1530.1111
$ws.Range("J102").Value = 8326.666999999999
$ws.Range("K102").Value = 1530.1111
$ws.Range("L102").Value = 8326.666999999999
$ws.Range("M102").Value = 91.88889999999992
$ws.Range("N102").Value = -11570.667
$ws.Range("H119").Value = 26500
$ws.Range("J119").Value = 26500
$ws.Range("L119").Value = 26500
$ws.Range("N119").Value = -36176

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2075.7368
$ws.Range("J86").Value = 2616.6667
$ws.Range("L86").Value = 2616.6667
$ws.Range("N86").Value = -4862.6667
$ws.Range("H89").Value = 2075.7368
$ws.Range("J89").Value = 2616.6667
$ws.Range("L89").Value = 13083.3335
$ws.Range("N89").Value = -24315.3335
$ws.Range("H99").Value = 1703.92
$ws.Range("I99").Value = 1462.2858
$ws.Range("J99").Value = 2972.5
$ws.Range("K99").Value = 1462.2858
$ws.Range("L99").Value = 2972.5
$ws.Range("M99").Value = 35.71419999999989
$ws.Range("N99").Value = -5968.5
$ws.Range("H134").Value = 4110.273
$ws.Range("I134").Value = 4311.35
$ws.Range("K134").Value = 12934.05
$ws.Range("M134").Value = -10399.05

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10636.548
$ws.Range("I31").Value = 12852.6875
$ws.Range("J31").Value = 3544.9
$ws.Range("K31").Value = 12852.6875
$ws.Range("L31").Value = 3544.9
$ws.Range("M31").Value = -12557.6875
$ws.Range("N31").Value = -4134.9
$ws.Range("H34").Value = 10636.548
$ws.Range("I34").Value = 12852.6875
$ws.Range("J34").Value = 3544.9
$ws.Range("K34").Value = 12852.6875
$ws.Range("L34").Value = 3544.9
$ws.Range("M34").Value = -12650.6875
$ws.Range("N34").Value = -3948.9
$ws.Range("H122").Value = 1291.9166
$ws.Range("I122").Value = 1148.3334
$ws.Range("J122").Value = 1435.5
$ws.Range("K122").Value = 3445.0002
$ws.Range("L122").Value = 4306.5
$ws.Range("M122").Value = -995.0001999999999
$ws.Range("N122").Value = -9206.5
$ws.Range("H134").Value = 1241.6364
$ws.Range("I134").Value = 1034
$ws.Range("J134").Value = 1449.2727
$ws.Range("K134").Value = 3102
$ws.Range("L134").Value = 4347.8181
$ws.Range("M134").Value = -567
$ws.Range("N134").Value = -9417.8181

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 130674.41
$ws.Range("J131").Value = 137790.39
$ws.Range("L131").Value = 413371.17
$ws.Range("N131").Value = -423451.17
$ws.Range("N130").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 22869.691
$ws.Range("I132").Value = 4304.4375
$ws.Range("J132").Value = 52574.1
$ws.Range("K132").Value = 12913.3125
$ws.Range("L132").Value = 157722.3
$ws.Range("M132").Value = -10383.3125
$ws.Range("N132").Value = -162782.3
$ws.Range("H135").Value = 49366
$ws.Range("J135").Value = 49184.445
$ws.Range("L135").Value = 49184.445
$ws.Range("N135").Value = -59324.445

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4518.8
$ws.Range("I61").Value = 2666.5
$ws.Range("J61").Value = 5753.6665
$ws.Range("K61").Value = 2666.5
$ws.Range("L61").Value = 5753.6665
$ws.Range("M61").Value = -2464.5
$ws.Range("N61").Value = -6157.6665
$ws.Range("H113").Value = 4518.8
$ws.Range("I113").Value = 2666.5
$ws.Range("J113").Value = 5753.6665
$ws.Range("K113").Value = 2666.5
$ws.Range("L113").Value = 5753.6665
$ws.Range("M113").Value = -496.5
$ws.Range("N113").Value = -10093.6665
$ws.Range("H128").Value = 35000
$ws.Range("J128").Value = 35000
$ws.Range("L128").Value = 35000
$ws.Range("H136").Value = 3076.2354
$ws.Range("I136").Value = 2176.6155
$ws.Range("K136").Value = 6529.8465
$ws.Range("M136").Value = -3979.8465
$ws.Range("N128").Value = -44960

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3247972
$ws.Range("I107").Value = 360.83334
$ws.Range("J107").Value = 5683680.5
$ws.Range("K107").Value = 1082.50002
$ws.Range("L107").Value = 17051041.5
$ws.Range("M107").Value = 837.4999800000001
$ws.Range("N107").Value = -17054881.5
$ws.Range("H113").Value = 1802693
$ws.Range("I113").Value = 1063.9166
$ws.Range("J113").Value = 9009209
$ws.Range("K113").Value = 3191.7498
$ws.Range("L113").Value = 27027627
$ws.Range("M113").Value = -1021.7498
$ws.Range("N113").Value = -27031967
$ws.Range("H126").Value = 1320.5
$ws.Range("I126").Value = 790.5833
$ws.Range("K126").Value = 2371.7499
$ws.Range("M126").Value = 98.2501000000002
